# Trade #77 closed at 2026-02-17 15:49:20 - unknown UNKNOWN +0.000%
#
# Updates the Summary and Strategy Status roll-up figures for the
# MarketMaking strategy and appends the newly closed trade (#77) as
# row 78 on both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.9    # Current Capital
$summary.Range("B4").Value = -0.11     # Total P&L $
$summary.Range("B5").Value = -0.03     # Total P&L %
$summary.Range("B6").Value = 77        # Total Trades
$summary.Range("B7").Value = 25        # Winning Trades
$summary.Range("B9").Value = 32.47     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (row 4 = MarketMaking)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.9       # Capital
$status.Range("D4").Value = 77         # Trades
$status.Range("E4").Value = -0.11      # P&L $
$status.Range("F4").Value = -0.1       # P&L %
$status.Range("G4").Value = 32.47      # Win Rate %

# ---------------------------------------------------------------------
# New trade row (shared by "All Trades" and "MarketMaking" sheets)
# ---------------------------------------------------------------------
function Add-TradeRow77($ws) {
    $row = 78

    # Force text formatting on the date/time columns so that Excel
    # does not reinterpret the strings as date/time serial numbers.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("C$row").NumberFormat = "@"

    $ws.Range("A$row").Value = 77
    $ws.Range("B$row").Value = "2026-02-17"
    $ws.Range("C$row").Value = "15:49:13"
    $ws.Range("D$row").Value = "MarketMaking"
    $ws.Range("E$row").Value = "UP"
    $ws.Range("F$row").Value = 0.89
    $ws.Range("G$row").Value = 0.97
    $ws.Range("H$row").Value = "CLOSED"
    $ws.Range("I$row").Value = 8.9888
    $ws.Range("J$row").Value = 0.08
    $ws.Range("K$row").Value = 99.9
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$row").Value = "early_exit"
    $ws.Range("Q$row").Value = 0.15
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow77 $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow77 $marketMaking
